# Updates to read SkipFirstRow from ItamaramExcelConfig.json and skip and
# read first row of excel as column headers.
#
# Data edits made to Book1.xlsx:
#   - Shared strings "great"/"teest"/"Threee" renamed to "Three"/"Four"/"Five"
#   - Row 1 (C1:E1) re-ordered/re-pointed so the visible header text becomes
#     C1="Four", D1="Five", E1="Three"
#   - Active selection moved from E1 to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared string values (affects every cell that references them).
# Order matters for how the shared-string table gets rebuilt on save, so
# update E1 first, then C1, then D1.
$ws.Cells.Item(1, 5).Value = "Three"
$ws.Cells.Item(1, 3).Value = "Four"
$ws.Cells.Item(1, 4).Value = "Five"

# Move the active selection to C2.
$ws.Range("C2").Select()
